# Commit: "Added new 350 node result"
# Update the Synthetic-350 / Uber-350 result rows (L10:U10 and L11:U11) with
# the newly-measured competitive-ratio values, and repoint the "350 node"
# chart's series-name references back at the shared A4/A5 label cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "Online/Offline" results for the 350-node synthetic run ---
$row10 = @(1.7432039164648001, 1.4856809842563701, 1.45240555106147, 1.46538313485802, 1.47191711977872, 1.4482665388335501, 1.44478059256847, 1.4577546370241901, 1.44853351425225, 1.4174839555005301)
$cols = @("L","M","N","O","P","Q","R","S","T","U")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])10").Value = $row10[$i]
}

# --- Row 11: "Greedy/Offline" results for the 350-node synthetic run ---
# (new run settled on a single constant competitive ratio across all coefficients)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])11").Value = 1.42981194919286
}

# --- Chart 3 (the "Synthetic Dataset 350 Nodes" chart) series naming ---
# Point the series names at the shared label cells (A4/A5) instead of the
# local ones (A10/A11), matching the other two charts' convention.
$co = $ws.ChartObjects(3)
$chart = $co.Chart
$series1 = $chart.SeriesCollection(1)
$series2 = $chart.SeriesCollection(2)
$series1.Name = '=Sheet1!$A$4'
$series2.Name = '=Sheet1!$A$5'
